$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.101275654079501
$ws.Cells.Item(2, 3).Value = 0.08445188921687219
$ws.Cells.Item(2, 4).Value = 0.02776715460236545
$ws.Cells.Item(2, 5).Value = 0.04050247491356451
$ws.Cells.Item(2, 6).Value = 1.748826202551683
$ws.Cells.Item(2, 7).Value = 1.498820550819858
$ws.Cells.Item(2, 8).Value = 0.02959398662645474
$ws.Cells.Item(2, 9).Value = 0.04258766848001017
$ws.Cells.Item(2, 10).Value = 0.9925520031292905
$ws.Cells.Item(2, 12).Value = 0.07309430288443863
$ws.Cells.Item(2, 13).Value = 0.8969209204041704
$ws.Cells.Item(2, 14).Value = 0.1524290226101002
$ws.Cells.Item(3, 2).Value = 0.9627377089899483
$ws.Cells.Item(3, 3).Value = 0.07673188251621355
$ws.Cells.Item(3, 4).Value = 0.02431779593750605
$ws.Cells.Item(3, 5).Value = 0.03674313781975869
$ws.Cells.Item(3, 6).Value = 1.661533148389637
$ws.Cells.Item(3, 7).Value = 1.421143582568547
$ws.Cells.Item(3, 8).Value = 0.03523965354584835
$ws.Cells.Item(3, 9).Value = 0.04956001605307403
$ws.Cells.Item(3, 10).Value = 0.9598053088124772
$ws.Cells.Item(3, 12).Value = 0.0675281700898509
$ws.Cells.Item(3, 13).Value = 0.7804142449179494
$ws.Cells.Item(3, 14).Value = 0.1341272434672831
$ws.Cells.Item(4, 2).Value = 0.8776866243974268
$ws.Cells.Item(4, 3).Value = 0.07205637888142036
$ws.Cells.Item(4, 4).Value = 0.02223026700075792
$ws.Cells.Item(4, 5).Value = 0.03442614592107951
$ws.Cells.Item(4, 6).Value = 1.608003426558042
$ws.Cells.Item(4, 7).Value = 1.373422510632423
$ws.Cells.Item(4, 8).Value = 0.03908550389980281
$ws.Cells.Item(4, 9).Value = 0.05430023649454885
$ws.Cells.Item(4, 10).Value = 0.939746154486869
$ws.Cells.Item(4, 12).Value = 0.06405456378271879
$ws.Cells.Item(4, 13).Value = 0.7091717536941644
$ws.Cells.Item(4, 14).Value = 0.1229808406250541
$ws.Cells.Item(5, 2).Value = 0.8426884463302144
$ws.Cells.Item(5, 3).Value = 0.07031912395294171
$ws.Cells.Item(5, 4).Value = 0.02145994430847509
$ws.Cells.Item(5, 5).Value = 0.03346095268117821
$ws.Cells.Item(5, 6).Value = 1.584059182156579
$ws.Cells.Item(5, 7).Value = 1.3517535142052
$ws.Cells.Item(5, 8).Value = 0.04075804107518199
$ws.Cells.Item(5, 9).Value = 0.05646021485543518
$ws.Cells.Item(5, 10).Value = 0.9304500677435072
$ws.Cells.Item(5, 12).Value = 0.06252926832616623
$ws.Cells.Item(5, 13).Value = 0.6805692058259183
$ws.Cells.Item(5, 14).Value = 0.1186525085266013
$ws.Cells.Item(6, 2).Value = 0.8364671859047519
$ws.Cells.Item(6, 3).Value = 0.07022195725193825
$ws.Cells.Item(6, 4).Value = 0.02142223322842796
$ws.Cells.Item(6, 5).Value = 0.03327807767114699
$ws.Cells.Item(6, 6).Value = 1.577468181413707
$ws.Cells.Item(6, 7).Value = 1.345456170783152
$ws.Cells.Item(6, 8).Value = 0.04105708517881923
$ws.Cells.Item(6, 9).Value = 0.05697481376816071
$ws.Cells.Item(6, 10).Value = 0.9275253586898771
$ws.Cells.Item(6, 12).Value = 0.06215900123655871
$ws.Cells.Item(6, 13).Value = 0.676271386999133
$ws.Cells.Item(6, 14).Value = 0.1181717247738661
$ws.Cells.Item(7, 2).Value = 0.8760907036875096
$ws.Cells.Item(7, 3).Value = 0.0725534283131779
$ws.Cells.Item(7, 4).Value = 0.02246528629283517
$ws.Cells.Item(7, 5).Value = 0.03435178009651274
$ws.Cells.Item(7, 6).Value = 1.600536675474316
$ws.Cells.Item(7, 7).Value = 1.365758430178872
$ws.Cells.Item(7, 8).Value = 0.03914917453510713
$ws.Cells.Item(7, 9).Value = 0.05471646700445909
$ws.Cells.Item(7, 10).Value = 0.9358497618637358
$ws.Cells.Item(7, 12).Value = 0.06371633851092895
$ws.Cells.Item(7, 13).Value = 0.7100102404061204
$ws.Cells.Item(7, 14).Value = 0.1235700878555264
$ws.Cells.Item(8, 2).Value = 1.051987155789988
$ws.Cells.Item(8, 3).Value = 0.08247070600808115
$ws.Cells.Item(8, 4).Value = 0.02689814062480167
$ws.Cells.Item(8, 5).Value = 0.0391261677804291
$ws.Cells.Item(8, 6).Value = 1.709267911693601
$ws.Cells.Item(8, 7).Value = 1.462291399072143
$ws.Cells.Item(8, 8).Value = 0.03150709736033797
$ws.Cells.Item(8, 9).Value = 0.04535642321617228
$ws.Cells.Item(8, 10).Value = 0.9762776784317424
$ws.Cells.Item(8, 12).Value = 0.07076392335160264
$ws.Cells.Item(8, 13).Value = 0.8583061285230826
$ws.Cells.Item(8, 14).Value = 0.146961955946594
$ws.Cells.Item(9, 2).Value = 1.399463102931037
$ws.Cells.Item(9, 3).Value = 0.1016132998264681
$ws.Cells.Item(9, 4).Value = 0.03541062003513673
$ws.Cells.Item(9, 5).Value = 0.04855217261871303
$ws.Cells.Item(9, 6).Value = 1.936164939233748
$ws.Cells.Item(9, 7).Value = 1.66490207322687
$ws.Cells.Item(9, 8).Value = 0.0196151336590602
$ws.Cells.Item(9, 9).Value = 0.0302215029166808
$ws.Cells.Item(9, 10).Value = 1.062795076450158
$ws.Cells.Item(9, 12).Value = 0.08474290489696301
$ws.Cells.Item(9, 13).Value = 1.149741250789702
$ws.Cells.Item(9, 14).Value = 0.1924846199784866
$ws.Cells.Item(10, 2).Value = 1.652215359775568
$ws.Cells.Item(10, 3).Value = 0.1176179022985124
$ws.Cells.Item(10, 4).Value = 0.04171561838003868
$ws.Cells.Item(10, 5).Value = 0.05342410143376775
$ws.Cells.Item(10, 6).Value = 2.064445017548522
$ws.Cells.Item(10, 7).Value = 1.776808846663357
$ws.Cells.Item(10, 8).Value = 0.01350100683356636
$ws.Cells.Item(10, 9).Value = 0.02201332731723848
$ws.Cells.Item(10, 10).Value = 1.10792522299684
$ws.Cells.Item(10, 12).Value = 0.09085576446199894
$ws.Cells.Item(10, 13).Value = 1.36962402603487
$ws.Cells.Item(10, 14).Value = 0.2199098545953717
$ws.Cells.Item(11, 2).Value = 1.740225740057838
$ws.Cells.Item(11, 3).Value = 0.1373657920589579
$ws.Cells.Item(11, 4).Value = 0.04270116333593421
$ws.Cells.Item(11, 5).Value = 0.03997867590524251
$ws.Cells.Item(11, 6).Value = 1.811588011508647
$ws.Cells.Item(11, 7).Value = 1.534851627911081
$ws.Cells.Item(11, 8).Value = 0.03173424669546421
$ws.Cells.Item(11, 9).Value = 0.02074952000771102
$ws.Cells.Item(11, 10).Value = 0.9803298832929954
$ws.Cells.Item(11, 12).Value = 0.06676713146386959
$ws.Cells.Item(11, 13).Value = 1.496419155908569
$ws.Cells.Item(11, 14).Value = 0.176167712983478
$ws.Cells.Item(12, 2).Value = 1.762202521322394
$ws.Cells.Item(12, 3).Value = 0.1518745503302199
$ws.Cells.Item(12, 4).Value = 0.04174167065401235
$ws.Cells.Item(12, 5).Value = 0.0304055534008203
$ws.Cells.Item(12, 6).Value = 1.596334429924966
$ws.Cells.Item(12, 7).Value = 1.333306943825846
$ws.Cells.Item(12, 8).Value = 0.07077844871241012
$ws.Cells.Item(12, 9).Value = 0.02055279991425429
$ws.Cells.Item(12, 10).Value = 0.8760064485652208
$ws.Cells.Item(12, 12).Value = 0.05428294141207068
$ws.Cells.Item(12, 13).Value = 1.555911854026448
$ws.Cells.Item(12, 14).Value = 0.1370821672723253
$ws.Cells.Item(13, 2).Value = 1.732568457593629
$ws.Cells.Item(13, 3).Value = 0.1638361008656375
$ws.Cells.Item(13, 4).Value = 0.03962020149047873
$ws.Cells.Item(13, 5).Value = 0.02318500582889271
$ws.Cells.Item(13, 6).Value = 1.386725616771287
$ws.Cells.Item(13, 7).Value = 1.139890356982519
$ws.Cells.Item(13, 8).Value = 0.12756925529726
$ws.Cells.Item(13, 9).Value = 0.02164039723798794
$ws.Cells.Item(13, 10).Value = 0.7776910982249774
$ws.Cells.Item(13, 12).Value = 0.04943791473668213
$ws.Cells.Item(13, 13).Value = 1.568356659500864
$ws.Cells.Item(13, 14).Value = 0.1005915587176105
$ws.Cells.Item(14, 2).Value = 1.688205784757201
$ws.Cells.Item(14, 3).Value = 0.1714940717795486
$ws.Cells.Item(14, 4).Value = 0.03768762573425022
$ws.Cells.Item(14, 5).Value = 0.01958645517386826
$ws.Cells.Item(14, 6).Value = 1.243928072972651
$ws.Cells.Item(14, 7).Value = 1.009494980557022
$ws.Cells.Item(14, 8).Value = 0.1779473605335227
$ws.Cells.Item(14, 9).Value = 0.02312074160552413
$ws.Cells.Item(14, 10).Value = 0.7121354571658998
$ws.Cells.Item(14, 12).Value = 0.05062508627981543
$ws.Cells.Item(14, 13).Value = 1.556166331129077
$ws.Cells.Item(14, 14).Value = 0.07722896624430575
$ws.Cells.Item(15, 2).Value = 1.666069232973143
$ws.Cells.Item(15, 3).Value = 0.1728922715343373
$ws.Cells.Item(15, 4).Value = 0.03707457729583297
$ws.Cells.Item(15, 5).Value = 0.01886084324577153
$ws.Cells.Item(15, 6).Value = 1.206865266307503
$ws.Cells.Item(15, 7).Value = 0.9758592783856983
$ws.Cells.Item(15, 8).Value = 0.1909465789700704
$ws.Cells.Item(15, 9).Value = 0.02394589905444633
$ws.Cells.Item(15, 10).Value = 0.6956434128334195
$ws.Cells.Item(15, 12).Value = 0.0514274170220439
$ws.Cells.Item(15, 13).Value = 1.543525401209791
$ws.Cells.Item(15, 14).Value = 0.07168861542716343
$ws.Cells.Item(16, 2).Value = 1.562912045645419
$ws.Cells.Item(16, 3).Value = 0.1635111917774026
$ws.Cells.Item(16, 4).Value = 0.03497308246966213
$ws.Cells.Item(16, 5).Value = 0.01843860432788968
$ws.Cells.Item(16, 6).Value = 1.19385734531312
$ws.Cells.Item(16, 7).Value = 0.9661666243027014
$ws.Cells.Item(16, 8).Value = 0.1807676175884154
$ws.Cells.Item(16, 9).Value = 0.0273392700224937
$ws.Cells.Item(16, 10).Value = 0.6953467617698408
$ws.Cells.Item(16, 12).Value = 0.04988109192957069
$ws.Cells.Item(16, 13).Value = 1.445979848713989
$ws.Cells.Item(16, 14).Value = 0.06874030547406562
$ws.Cells.Item(17, 2).Value = 1.507450162811352
$ws.Cells.Item(17, 3).Value = 0.1521724886674747
$ws.Cells.Item(17, 4).Value = 0.0344007487012945
$ws.Cells.Item(17, 5).Value = 0.01973859696582769
$ws.Cells.Item(17, 6).Value = 1.257636031169113
$ws.Cells.Item(17, 7).Value = 1.025883774897835
$ws.Cells.Item(17, 8).Value = 0.1444424473821471
$ws.Cells.Item(17, 9).Value = 0.02922331215172669
$ws.Cells.Item(17, 10).Value = 0.7290030800967884
$ws.Cells.Item(17, 12).Value = 0.04723889893544264
$ws.Cells.Item(17, 13).Value = 1.377955767711114
$ws.Cells.Item(17, 14).Value = 0.0783128856416937
$ws.Cells.Item(18, 2).Value = 1.487205913432746
$ws.Cells.Item(18, 3).Value = 0.1380446057327589
$ws.Cells.Item(18, 4).Value = 0.0348498233695409
$ws.Cells.Item(18, 5).Value = 0.02404034873335004
$ws.Cells.Item(18, 6).Value = 1.407061130473878
$ws.Cells.Item(18, 7).Value = 1.164725553478476
$ws.Cells.Item(18, 8).Value = 0.09189941705393068
$ws.Cells.Item(18, 9).Value = 0.0294420509538087
$ws.Cells.Item(18, 10).Value = 0.801920446379782
$ws.Cells.Item(18, 12).Value = 0.04758714182761004
$ws.Cells.Item(18, 13).Value = 1.326553699937023
$ws.Cells.Item(18, 14).Value = 0.1018352556688811
$ws.Cells.Item(19, 2).Value = 1.494214480363354
$ws.Cells.Item(19, 3).Value = 0.1256460987376613
$ws.Cells.Item(19, 4).Value = 0.03652899888265182
$ws.Cells.Item(19, 5).Value = 0.03243859456328835
$ws.Cells.Item(19, 6).Value = 1.616989302746276
$ws.Cells.Item(19, 7).Value = 1.359713553378469
$ws.Cells.Item(19, 8).Value = 0.04560585722578736
$ws.Cells.Item(19, 9).Value = 0.02887012592672189
$ws.Cells.Item(19, 10).Value = 0.9020109341258973
$ws.Cells.Item(19, 12).Value = 0.05649123150509894
$ws.Cells.Item(19, 13).Value = 1.295046059936539
$ws.Cells.Item(19, 14).Value = 0.1396055225052635
$ws.Cells.Item(20, 2).Value = 1.581952260166304
$ws.Cells.Item(20, 3).Value = 0.1151507176029938
$ws.Cells.Item(20, 4).Value = 0.04082253574424755
$ws.Cells.Item(20, 5).Value = 0.05186715690589061
$ws.Cells.Item(20, 6).Value = 2.007194535487173
$ws.Cells.Item(20, 7).Value = 1.723272583824155
$ws.Cells.Item(20, 8).Value = 0.0150300135565562
$ws.Cells.Item(20, 9).Value = 0.02504553855298486
$ws.Cells.Item(20, 10).Value = 1.083804457945064
$ws.Cells.Item(20, 12).Value = 0.08812564893759856
$ws.Cells.Item(20, 13).Value = 1.315655363809725
$ws.Cells.Item(20, 14).Value = 0.2144313373635924
$ws.Cells.Item(21, 2).Value = 1.781500708600049
$ws.Cells.Item(21, 3).Value = 0.1255857523685648
$ws.Cells.Item(21, 4).Value = 0.045975907927577
$ws.Cells.Item(21, 5).Value = 0.05878857037064122
$ws.Cells.Item(21, 6).Value = 2.166507252664871
$ws.Cells.Item(21, 7).Value = 1.866747009971363
$ws.Cells.Item(21, 8).Value = 0.01031780356902593
$ws.Cells.Item(21, 9).Value = 0.01908477050680002
$ws.Cells.Item(21, 10).Value = 1.147635087405206
$ws.Cells.Item(21, 12).Value = 0.09881589679304881
$ws.Cells.Item(21, 13).Value = 1.480388061224147
$ws.Cells.Item(21, 14).Value = 0.2461213462567713
$ws.Cells.Item(22, 2).Value = 1.912409066634609
$ws.Cells.Item(22, 3).Value = 0.1324318566899336
$ws.Cells.Item(22, 4).Value = 0.04892432322314022
$ws.Cells.Item(22, 5).Value = 0.06237615653942186
$ws.Cells.Item(22, 6).Value = 2.261951715108935
$ws.Cells.Item(22, 7).Value = 1.953111904052605
$ws.Cells.Item(22, 8).Value = 0.007911748635897942
$ws.Cells.Item(22, 9).Value = 0.01543120032334233
$ws.Cells.Item(22, 10).Value = 1.185642439989664
$ws.Cells.Item(22, 12).Value = 0.1042843682956018
$ws.Cells.Item(22, 13).Value = 1.588533492716948
$ws.Cells.Item(22, 14).Value = 0.262558012330075
$ws.Cells.Item(23, 2).Value = 1.843922574320544
$ws.Cells.Item(23, 3).Value = 0.1281097926415384
$ws.Cells.Item(23, 4).Value = 0.04704644524230872
$ws.Cells.Item(23, 5).Value = 0.06053438017566037
$ws.Cells.Item(23, 6).Value = 2.21932783179389
$ws.Cells.Item(23, 7).Value = 1.915615981220441
$ws.Cells.Item(23, 8).Value = 0.009132220702440352
$ws.Cells.Item(23, 9).Value = 0.01694530627816082
$ws.Cells.Item(23, 10).Value = 1.169708152302036
$ws.Cells.Item(23, 12).Value = 0.1017485461909367
$ws.Cells.Item(23, 13).Value = 1.529293354459412
$ws.Cells.Item(23, 14).Value = 0.2529790533678238
$ws.Cells.Item(24, 2).Value = 1.583095977629455
$ws.Cells.Item(24, 3).Value = 0.1130831459615251
$ws.Cells.Item(24, 4).Value = 0.0404902941468066
$ws.Cells.Item(24, 5).Value = 0.05346868770591762
$ws.Cells.Item(24, 6).Value = 2.045809311236908
$ws.Cells.Item(24, 7).Value = 1.760885288595659
$ws.Cells.Item(24, 8).Value = 0.01470492198159751
$ws.Cells.Item(24, 9).Value = 0.0243592919792528
$ws.Cells.Item(24, 10).Value = 1.102982689291878
$ws.Cells.Item(24, 12).Value = 0.09151014693204118
$ws.Cells.Item(24, 13).Value = 1.308602393260145
$ws.Cells.Item(24, 14).Value = 0.2182609711090606
$ws.Cells.Item(25, 2).Value = 1.303236828394631
$ws.Cells.Item(25, 3).Value = 0.09732165385729985
$ws.Cells.Item(25, 4).Value = 0.03352881453209022
$ws.Cells.Item(25, 5).Value = 0.04589157688120515
$ws.Cells.Item(25, 6).Value = 1.861773341097006
$ws.Cells.Item(25, 7).Value = 1.596749308855024
$ws.Cells.Item(25, 8).Value = 0.02252443318472186
$ws.Cells.Item(25, 9).Value = 0.03449564999337706
$ws.Cells.Item(25, 10).Value = 1.032562797439908
$ws.Cells.Item(25, 12).Value = 0.08042165109742427
$ws.Cells.Item(25, 13).Value = 1.072747458119579
$ws.Cells.Item(25, 14).Value = 0.1812604743500117
